$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Worksheet, [string]$Address, [string]$Text)
    $range = $Worksheet.Range($Address)
    $range.NumberFormat = "@"
    $range.Value = $Text
    $range.Style = "Normal"
}

Set-TextValue $ws "D2" "26.495.89"
Set-TextValue $ws "E2" "  -0.78%  "
Set-TextValue $ws "D3" "1.850.79"
Set-TextValue $ws "E3" "  -1.29%  "
Set-TextValue $ws "D4" "1.000"
Set-TextValue $ws "E4" "  +0.11%  "
Set-TextValue $ws "D5" "261.52"
Set-TextValue $ws "E5" "  -7.43%  "
Set-TextValue $ws "D6" "0.9999"
Set-TextValue $ws "E6" "  +0.10%  "
Set-TextValue $ws "D7" "0.5143"
Set-TextValue $ws "D8" "0.3272"
Set-TextValue $ws "E8" "  -7.36%  "
Set-TextValue $ws "D9" "0.06781"
Set-TextValue $ws "E9" "  -5.47%  "
Set-TextValue $ws "D10" "19.08"
Set-TextValue $ws "E10" "  -5.37%  "
Set-TextValue $ws "D11" "0.7745"
Set-TextValue $ws "E11" "  -5.87%  "
Set-TextValue $ws "D12" "0.07701"
Set-TextValue $ws "E12" "  -0.63%  "
Set-TextValue $ws "D13" "1.880.30"
Set-TextValue $ws "E13" "  +0.32%  "
Set-TextValue $ws "D14" "88.74"
Set-TextValue $ws "E14" "  -0.98%  "
Set-TextValue $ws "D15" "5.041"
Set-TextValue $ws "E15" "  -2.14%  "
Set-TextValue $ws "D16" "1.0000"
Set-TextValue $ws "E16" "  +0.13%  "
Set-TextValue $ws "D17" "14.15"
Set-TextValue $ws "E17" "  -1.93%  "
Set-TextValue $ws "D18" "0.9996"
Set-TextValue $ws "E18" "  +0.16%  "
Set-TextValue $ws "D19" "0.000007912"
Set-TextValue $ws "E19" "  -3.71%  "
Set-TextValue $ws "D20" "26.531.15"
Set-TextValue $ws "E20" "  -0.84%  "
Set-TextValue $ws "D21" "2.069.10"
Set-TextValue $ws "E21" "  -1.39%  "
Set-TextValue $ws "D22" "4.556"
Set-TextValue $ws "E22" "  -5.18%  "
Set-TextValue $ws "D23" "9.552"
Set-TextValue $ws "E23" "  -5.97%  "
Set-TextValue $ws "D24" "5.954"
Set-TextValue $ws "E24" "  -4.79%  "
Set-TextValue $ws "D25" "2.368"
Set-TextValue $ws "E25" "  -1.64%  "
Set-TextValue $ws "D26" "144.58"
Set-TextValue $ws "E26" "  -0.67%  "
Set-TextValue $ws "D27" "1.659"
Set-TextValue $ws "E27" "  -0.59%  "
Set-TextValue $ws "D28" "16.96"
Set-TextValue $ws "E28" "  -2.65%  "
Set-TextValue $ws "D29" "111.36"
Set-TextValue $ws "E29" "  +0.39%  "
Set-TextValue $ws "D30" "4.220"
Set-TextValue $ws "E30" "  -4.43%  "
Set-TextValue $ws "D31" "4.169"
Set-TextValue $ws "E31" "  -4.51%  "
Set-TextValue $ws "D32" "0.08749"
Set-TextValue $ws "E32" "  -0.97%  "
Set-TextValue $ws "D33" "0.04854"
Set-TextValue $ws "E33" "  -1.18%  "
Set-TextValue $ws "D34" "1.138"
Set-TextValue $ws "E34" "  -3.46%  "
Set-TextValue $ws "D35" "2.850"
Set-TextValue $ws "E35" "  -0.47%  "
Set-TextValue $ws "D36" "0.6930"
Set-TextValue $ws "E36" "  -7.44%  "
Set-TextValue $ws "D37" "3.118"
Set-TextValue $ws "E37" "  -5.33%  "
Set-TextValue $ws "D38" "0.01809"
Set-TextValue $ws "E38" "  -3.81%  "
Set-TextValue $ws "D39" "2.225"
Set-TextValue $ws "E39" "  -8.42%  "
Set-TextValue $ws "D40" "0.4928"
Set-TextValue $ws "E40" "  -7.02%  "
Set-TextValue $ws "D41" "113.53"
Set-TextValue $ws "E41" "  -3.10%  "
Set-TextValue $ws "D42" "0.8969"
Set-TextValue $ws "E42" "  -7.90%  "
Set-TextValue $ws "D43" "6.158"
Set-TextValue $ws "E43" "  -2.36%  "
Set-TextValue $ws "D44" "7.811"
Set-TextValue $ws "E44" "  -4.83%  "
Set-TextValue $ws "D45" "0.9995"
Set-TextValue $ws "E45" "  +0.12%  "
Set-TextValue $ws "D46" "0.4266"
Set-TextValue $ws "E46" "  -7.66%  "
Set-TextValue $ws "D47" "0.1269"
Set-TextValue $ws "E47" "  -7.09%  "
Set-TextValue $ws "D48" "9.178"
Set-TextValue $ws "E48" "  -3.24%  "
Set-TextValue $ws "D49" "0.05899"
Set-TextValue $ws "E49" "  -0.48%  "
Set-TextValue $ws "E50" "  -4.38%  "
